$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I (roboticS1Prep) changes from text "No" to a boolean FALSE value,
# displayed via a custom TRUE/FALSE number format. Column J keeps its
# "random" text (its underlying shared-string index shifts down once "No"
# is dropped, but that's handled automatically by the engine).
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Value = $false
    $cell.NumberFormat = '"TRUE";"TRUE";"FALSE"'
}

# Widen column I (9) to fit the new boolean header/values; leave every
# other column at its existing default width.
$ws.Columns.Item(9).ColumnWidth = 15.5

# Move the active selection from column H to column I.
$ws.Range("I2:I27").Select()
